$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P needs the same header style (bold, centered, bordered) as the rest of row 1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null

# --- Update header row (H1:P1) ---
$ws.Range("H1").Value = 'Glucose2 [g/L]'
$ws.Range("I1").Value = 'Glucose [g/L].1'
$ws.Range("J1").Value = 'added Glucose  [g]'
$ws.Range("K1").Value = 'Consumed_Glucose [g]'
$ws.Range("L1").Value = 'Consumed_Glucose [g/L]'
$ws.Range("M1").Value = 'Volume at each time point'
$ws.Range("N1").Value = 'mu between direct samples [1/h]'
$ws.Range("O1").Value = 'Yxs [gx/gs]'
$ws.Range("P1").Value = 'qs 1 - ds/dt/x - [1/h]'

# --- Update data rows 2-21 ---
# row 2
$ws.Range("C2").Value = 10.44
$ws.Range("H2").Value = 5.22
$ws.Range("I2").Value = 5.22
$ws.Range("J2").Value = 10.44
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = $null
$ws.Range("O2").Value = $null
$ws.Range("P2").Value = $null

# row 3
$ws.Range("C3").Value = 10.352
$ws.Range("H3").Value = 5.176
$ws.Range("I3").Value = 5.176
$ws.Range("J3").Value = 10.44
$ws.Range("K3").Value = 0.0879999999999991
$ws.Range("L3").Value = 0.0879999999999991
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0.4676677286516002
$ws.Range("O3").Value = 7.007575761363708
$ws.Range("P3").Value = 0.06673744880934256

# row 4
$ws.Range("C4").Value = 10.396
$ws.Range("H4").Value = 5.198
$ws.Range("I4").Value = 5.198
$ws.Range("J4").Value = 10.44
$ws.Range("K4").Value = 0.0439999999999987
$ws.Range("L4").Value = 0.0439999999999987
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 0.8825025417389376
$ws.Range("O4").Value = 5.257575749999952
$ws.Range("P4").Value = 0.1678535096216057

# row 5
$ws.Range("C5").Value = 9.37
$ws.Range("H5").Value = 4.685
$ws.Range("I5").Value = 4.685
$ws.Range("J5").Value = 10.44
$ws.Range("K5").Value = 1.070000000000001
$ws.Range("L5").Value = 1.070000000000001
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 0.5408167273906378
$ws.Range("O5").Value = 0.8940870692007782
$ws.Range("P5").Value = 0.6048814998231351

# row 6
$ws.Range("C6").Value = 8.46
$ws.Range("H6").Value = 4.23
$ws.Range("I6").Value = 4.23
$ws.Range("J6").Value = 10.44
$ws.Range("K6").Value = 1.979999999999998
$ws.Range("L6").Value = 1.979999999999998
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 0.3160748173374999
$ws.Range("O6").Value = 0.854212454945057
$ws.Range("P6").Value = 0.370018975382219

# row 7
$ws.Range("C7").Value = 6.96
$ws.Range("H7").Value = 3.48
$ws.Range("I7").Value = 3.48
$ws.Range("J7").Value = 10.44
$ws.Range("K7").Value = 3.48
$ws.Range("L7").Value = 3.48
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 0.15458832444611
$ws.Range("O7").Value = 0.3168888886666664
$ws.Range("P7").Value = 0.4878313187204257

# row 8
$ws.Range("C8").Value = 5.762
$ws.Range("H8").Value = 2.881
$ws.Range("I8").Value = 2.881
$ws.Range("J8").Value = 10.44
$ws.Range("K8").Value = 4.678
$ws.Range("L8").Value = 4.678
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 0.2482040228987069
$ws.Range("O8").Value = 0.7690595434056761
$ws.Range("P8").Value = 0.3227370689655173

# row 9
$ws.Range("C9").Value = 4.036
$ws.Range("H9").Value = 2.018
$ws.Range("I9").Value = 2.018
$ws.Range("J9").Value = 10.44
$ws.Range("K9").Value = 6.404
$ws.Range("L9").Value = 6.404
$ws.Range("M9").Value = 1
$ws.Range("N9").Value = 0.08027262416329092
$ws.Range("O9").Value = 0.204712244495944
$ws.Range("P9").Value = 0.3921241954087479

# row 10
$ws.Range("C10").Value = 3.134
$ws.Range("H10").Value = 1.567
$ws.Range("I10").Value = 1.567
$ws.Range("J10").Value = 10.44
$ws.Range("K10").Value = 7.305999999999999
$ws.Range("L10").Value = 7.305999999999999
$ws.Range("M10").Value = 1
$ws.Range("N10").Value = 0.07259358288284798
$ws.Range("O10").Value = 0.4013303769401335
$ws.Range("P10").Value = 0.1808823529290851

# row 11
$ws.Range("C11").Value = 2.082
$ws.Range("H11").Value = 1.041
$ws.Range("I11").Value = 1.041
$ws.Range("J11").Value = 10.44
$ws.Range("K11").Value = 8.358
$ws.Range("L11").Value = 8.358
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 0.03789106306631612
$ws.Range("O11").Value = 0.1926489220532315
$ws.Range("P11").Value = 0.1966845319583273

# row 12
$ws.Range("C12").Value = 56.204
$ws.Range("H12").Value = 28.102
$ws.Range("I12").Value = 28.102
$ws.Range("J12").Value = 145.862
$ws.Range("K12").Value = 89.65799999999999
$ws.Range("L12").Value = 68.74173270281385
$ws.Range("M12").Value = 1.304273204570096
$ws.Range("N12").Value = 0.04532945039001619
$ws.Range("O12").Value = 0.0500079496883982
$ws.Range("P12").Value = 0.9064448887120156

# row 13
$ws.Range("C13").Value = 69.7
$ws.Range("H13").Value = 34.85
$ws.Range("I13").Value = 34.85
$ws.Range("J13").Value = 169.662
$ws.Range("K13").Value = 99.962
$ws.Range("L13").Value = 71.68328608621864
$ws.Range("M13").Value = 1.394495222774365
$ws.Range("N13").Value = 0.03949364134873408
$ws.Range("O13").Value = 0.2301505061303305
$ws.Range("P13").Value = 0.1715991939916459

# row 14
$ws.Range("C14").Value = 79.256
$ws.Range("H14").Value = 39.628
$ws.Range("I14").Value = 39.628
$ws.Range("J14").Value = 181.562
$ws.Range("K14").Value = 102.306
$ws.Range("L14").Value = 70.83332697721671
$ws.Range("M14").Value = 1.444320129603772
$ws.Range("N14").Value = 0.06718569777249125
$ws.Range("O14").Value = 0.7310155587715322
$ws.Range("P14").Value = 0.09190734310141971

# row 15
$ws.Range("C15").Value = 94.756
$ws.Range("H15").Value = 47.378
$ws.Range("I15").Value = 47.378
$ws.Range("J15").Value = 205.362
$ws.Range("K15").Value = 110.606
$ws.Range("L15").Value = 71.18688537354446
$ws.Range("M15").Value = 1.553741246292889
$ws.Range("N15").Value = 0.01931993821409123
$ws.Range("O15").Value = 0.9427956950314005
$ws.Range("P15").Value = 0.02049217907539105

# row 16
$ws.Range("C16").Value = 104.982
$ws.Range("H16").Value = 52.491
$ws.Range("I16").Value = 52.491
$ws.Range("J16").Value = 217.262
$ws.Range("K16").Value = 112.28
$ws.Range("L16").Value = 69.55804999721735
$ws.Range("M16").Value = 1.614191312213205
$ws.Range("N16").Value = 0.006430868127268053
$ws.Range("O16").Value = 0.03274323100733601
$ws.Range("P16").Value = 0.1964029794685576

# row 17
$ws.Range("C17").Value = 113.33
$ws.Range("H17").Value = 56.665
$ws.Range("I17").Value = 56.665
$ws.Range("J17").Value = 224.402
$ws.Range("K17").Value = 111.072
$ws.Range("L17").Value = 67.20575159889712
$ws.Range("M17").Value = 1.652715688129031
$ws.Range("N17").Value = 0.1912081984223309
$ws.Range("O17").Value = 0.4018764514209001
$ws.Range("P17").Value = 0.4757885109871031

# row 18
$ws.Range("C18").Value = 119.556
$ws.Range("H18").Value = 59.778
$ws.Range("I18").Value = 59.778
$ws.Range("J18").Value = 236.3
$ws.Range("K18").Value = 116.744
$ws.Range("L18").Value = 69.0470995585246
$ws.Range("M18").Value = 1.690787893285037
$ws.Range("N18").Value = 0.02866889234753486
$ws.Range("O18").Value = 0.05720446939388425
$ws.Range("P18").Value = 0.5011652524933631

# row 19
$ws.Range("C19").Value = 136.44
$ws.Range("H19").Value = 150.34
$ws.Range("I19").Value = 68.22
$ws.Range("J19").Value = 314.84
$ws.Range("K19").Value = 178.4
$ws.Range("L19").Value = 90.54990952549136
$ws.Range("M19").Value = 1.970184188309733
$ws.Range("N19").Value = 0.1700707514016821
$ws.Range("O19").Value = 0.1579948544036376
$ws.Range("P19").Value = 1.076432217008749

# row 20
$ws.Range("C20").Value = 141.28
$ws.Range("H20").Value = 170.96
$ws.Range("I20").Value = 70.64
$ws.Range("J20").Value = 350.84
$ws.Range("K20").Value = 209.56
$ws.Range("L20").Value = 99.5918914827856
$ws.Range("M20").Value = 2.104187367866412
$ws.Range("N20").Value = 0.02479162268240855
$ws.Range("O20").Value = 0.03421078934474742
$ws.Range("P20").Value = 0.724672629812178

# row 21
$ws.Range("C21").Value = 120.475
$ws.Range("H21").Value = 120.475
$ws.Range("I21").Value = $null
$ws.Range("J21").Value = 360.44
$ws.Range("K21").Value = 239.965
$ws.Range("L21").Value = 98.09957483368598
$ws.Range("M21").Value = 2.446137003211552
$ws.Range("N21").Value = 0.0477680122771348
$ws.Range("O21").Value = 4.297568705589018
$ws.Range("P21").Value = 0.01111512474832855

